# Corrections on the MAIN sheet: the "BUSINESS FLOW" column (D) used to hold
# pipe-delimited header strings describing the fields involved in each test
# step. They are replaced with the actual automation step/function names
# used by the article search, and a new step (row 4 / T3) is added for the
# "complete order" flow.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAIN")

$ws.Range("D2").Value = "login"
$ws.Range("D3").Value = "addItem"
$ws.Range("D4").Value = "completeOrder"

# Move/restore the active selection to D4 (was parked on D13, an empty cell
# below the used range).
$ws.Range("D4").Select() | Out-Null
